$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 37 ---------------------------------------------------------------
$ws.Range("A37").Value = "2023-03-10-1326_RF_avg_restricted.csv"
$ws.Range("B37").Value = "RandomForest"
$ws.Range("C37").Value = "MoCo"
$ws.Range("D37").Value = "Centers"
$ws.Range("E37").Value = "1 x 3"
$ws.Range("F37").Value = "average"
$ws.Range("H37").Value = "{""max_depth"": 3}"
$ws.Range("I37").Value = "submission trained on all data"
$ws.Range("N37").Value = "March 10, 2023, 12:27 p.m."
$ws.Range("O37").Value = 0.607

# --- Row 38 ----------------------------------------------------------------
$ws.Range("B38").Value = "RandomForest"
$ws.Range("C38").Value = "MoCo"
$ws.Range("D38").Value = "Centers"
$ws.Range("E38").Value = "1 x 3"
$ws.Range("F38").Value = "average"
$ws.Range("H38").Value = "{""max_depth"": 3}"
$ws.Range("J38").Value = 0.583
$ws.Range("K38").Value = 0.734
$ws.Range("L38").Value = 0.598
$ws.Range("M38").Value = "0.639 (0.068)"
$ws.Range("N38").Value = "March 10, 2023, 12:34 p.m."
$ws.Range("O38").Value = 0.603
$ws.Range("A38").Value = "2023-03-10-1333_RF_avg_restricted_cv.csv"

# Highlight the hyperparameters cell of the new cross-validated submission
# with a thin black top/bottom rule, like the other "restricted" entries.
$h38 = $ws.Range("H38")
$h38.Font.Color = 0
$h38.Borders.Item(8).Color = 0
$h38.Borders.Item(9).Color = 0

# --- Extend the table + used range to cover the two new rows ---------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:O38"))

# --- Update the view so the new rows are in focus ---------------------------
$ws.Range("A39").Select()
